$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (matches original t="str" cells) for D and G columns in data rows 8-19
$ws.Range("D8:D19").NumberFormat = "@"
$ws.Range("G8:G19").NumberFormat = "@"

# Row 8
$ws.Range("A8").Value = 'P. point'
$ws.Range("C8").Value = 57
$ws.Range("D8").Value = '2'
$ws.Range("E8").Value = 'Short point (up to 3 mtr.)'
$ws.Range("F8").Value = 256
$ws.Range("G8").Value = '14592.00'

# Row 9
$ws.Range("A9").Value = 'P. point'
$ws.Range("C9").Value = 4
$ws.Range("D9").Value = '3'
$ws.Range("E9").Value = 'Medium point (up to 6 mtr.)'
$ws.Range("F9").Value = 472
$ws.Range("G9").Value = '1888.00'

# Row 10
$ws.Range("A10").Value = ''
$ws.Range("C10").Value = 23
$ws.Range("D10").Value = '2.0'
$ws.Range("E10").Value = 'Rewiring of 3/5 pin 6 amp. Light plug point with 1.5 sq. mm nominal size  FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade(IS:694)   in recessed ISI marked MMS ( IS:9537 P - III ) virgin material  PVC conduit & it''s  ISI marked (IS:3419-1988) accessories, 1.2 mm thick  MS box with earth terminal of required size,  6 A  switch, 3/5 pin 6 A socket, 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/ brass  screws, cup washers, making connections, testing etc. as required.  For specification of copper  Conductor,  Phenolic Laminated sheet''s & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = '0.00'

# Row 11
$ws.Range("C11").Value = 24
$ws.Range("D11").Value = '3.0'
$ws.Range("E11").Value = 'P & F ISI marked (IS:3854) 6 amp. flush type non modular switch  with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F11").Value = 23
$ws.Range("G11").Value = '552.00'

# Row 12
$ws.Range("A12").Value = 'Each'
$ws.Range("C12").Value = 94
$ws.Range("D12").Value = '9.0'
$ws.Range("E12").Value = 'Providing & Fixing of IS 11037:1984  marked  non modular socket size flush type 180 watt rotary minimum 5 step fan regulator with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F12").Value = 219
$ws.Range("G12").Value = '20586.00'

# Row 13
$ws.Range("C13").Value = 93
$ws.Range("G13").Value = '5208.00'

# Row 14
$ws.Range("A14").Value = 'Set'
$ws.Range("C14").Value = 97
$ws.Range("D14").Value = '13.0'
$ws.Range("E14").Value = 'Plate Earthing  as per IS:3043 with Hot dipped G.I. Earth plate of size 600mm x 600mm x 6.0mm by embodying  3 to 4 mtr. below the ground level with 20  mm dia. G.I. ''B'' class watering Pipe ,including all accessories like nut, bolts, reducer, nipple, wire meshed funnel, and Heavy duty weather proof poly-propylene earth pit chamber with lockable Jam free lid suitable for safe working load 5000 Kg or more of size Top Dia. 225 to 260 mm, Bottom Dia 300 to 350 mm. and Height  250 to 300 mm. and embodying the pipe  complete with alternate layers salt and coke/ charcoal, testing of earth resistance for value of 5 ohms or less  as required & must record by engineer in charge during site visit and ensure to enter in measurment book.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .   '
$ws.Range("F14").Value = 5733
$ws.Range("G14").Value = '556101.00'

# Row 15
$ws.Range("A15").Value = ''
$ws.Range("C15").Value = 30
$ws.Range("D15").Value = '15.0'
$ws.Range("E15").Value = 'Providing & Fixing of  BEE  Star rated copper wounded double ball bearing capacitor start, aluminium body & Metallic  blade ceiling  fan  Conforming to all the performance requirements laid down in IS 374:2019 including all amendments, as applicable ; & Carry BIS licensing (i.e. ISI marking) with down rod up to 80 cm with secondary support safety cable ( steel rope) , cotter pin with 3 x 1.5 sq.mm pvc insulated flexible copper conductor making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = '0.00'

# Row 16
$ws.Range("A16").Value = 'Each'
$ws.Range("C16").Value = 23
$ws.Range("D16").Value = '27'
$ws.Range("E16").Value = '1170mm(+/-10%) LED batten with min. lumen output 2200 lm'
$ws.Range("F16").Value = 492
$ws.Range("G16").Value = '11316.00'

# Row 17
$ws.Range("A17").Value = ''
$ws.Range("C17").Value = 66
$ws.Range("D17").Value = '18.0'
$ws.Range("E17").Value = 'Providing & Fixing of Recessed/surface mounting heavy duty horizontal type Double Door ( Metal / Glazed )Distribution board with Metal end box made out from Galvanized steel / CRCA sheet not less then 1.2 mm thick  conforming to IS-8623-1 & 3 /  IEC 61439- 1 & 3, powder painted complete with reversible door (for double door DB only )100 amp.  insulated copper bus bar/shorting link , copper neutral link, copper earth link , color coded interconnecting wire set  of suitable rating and din bar,masking sheet,  making internal DB  terminations with copper lugs, Ferrules,  detachable gland plate, including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = '0.00'

# Row 18
$ws.Range("A18").Value = '%'
$ws.Range("C18").Value = 17
$ws.Range("D18").Value = '37'
$ws.Range("E18").Value = 'Add Tender Premium '

# Row 19
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = '38'
$ws.Range("E19").Value = 'Grand Total'

# Remove the old row 20 ("Grand Total" calc row) entirely; rows below shift up by one,
# turning old row 21 (spacer) into new row 20, and old rows 22-24 into new rows 21-23.
$ws.Rows(20).Delete()

# Re-apply text storage for the summary amount cells, then set the refreshed totals
$ws.Range("G21").NumberFormat = "@"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("H23").NumberFormat = "@"

$ws.Range("G21").Value = '610243.00'
$ws.Range("H21").Value = '610243.00'
$ws.Range("G23").Value = '610243.00'
$ws.Range("H23").Value = '610243.00'
